# Mark the "Add CheatSystem to the engine" and "Remove the existing,
# copy-and-paste cheat code" bullets as done by giving them the same
# green highlight already used on the "Implemented the CheatSystem
# functionality" bullet above them (slide 22, "Content Placeholder 2").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(22)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$green = 65280  # RGB(0, 255, 0) -> srgbClr val="00FF00"

foreach ($text in @("Add CheatSystem to the engine", "Remove the existing, copy-and-paste cheat code")) {
    $count = $tr.Paragraphs().Count
    for ($i = 1; $i -le $count; $i++) {
        $para = $tr.Paragraphs($i, 1)
        if ($para.Text.TrimEnd() -eq $text) {
            $para.Font.Highlight = $green
            break
        }
    }
}
